$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Crime statistics table updates (rows 15-27) ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = -28.571428571428
$ws.Range("N15").Value = -37.5
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 57
$ws.Range("K16").Value = 7.017543859649
$ws.Range("L16").Value = 29.787234042553
$ws.Range("M16").Value = -18.666666666666
$ws.Range("N16").Value = -59.333333333333
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 75
$ws.Range("I17").Value = 94
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 28.767123287671
$ws.Range("L17").Value = 40.298507462686
$ws.Range("M17").Value = 59.322033898305
$ws.Range("N17").Value = 27.027027027027
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = 36.111111111111
$ws.Range("L18").Value = 32.432432432432
$ws.Range("M18").Value = -38.75
$ws.Range("N18").Value = -82.926829268292
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 87.5
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 47.058823529411
$ws.Range("I19").Value = 211
$ws.Range("J19").Value = 182
$ws.Range("K19").Value = 15.934065934065
$ws.Range("L19").Value = 63.565891472868
$ws.Range("M19").Value = 63.565891472868
$ws.Range("N19").Value = 57.462686567164
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 650
$ws.Range("F20").Value = 48
$ws.Range("H20").Value = 128.571428571429
$ws.Range("I20").Value = 171
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = 103.571428571429
$ws.Range("L20").Value = 216.666666666667
$ws.Range("M20").Value = 222.641509433962
$ws.Range("N20").Value = -72.813990461049
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 80.952380952380
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 62.5
$ws.Range("I21").Value = 592
$ws.Range("J21").Value = 437
$ws.Range("K21").Value = 35.469107551487
$ws.Range("L21").Value = 74.631268436578
$ws.Range("M21").Value = 46.898263027295
$ws.Range("N21").Value = -53.929961089494
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("L22").Value = 75
$ws.Range("M22").Value = 16.666666666666
$ws.Range("C23").Value = 1
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 12
$ws.Range("K23").Value = 33.333333333333
$ws.Range("L23").Value = 45.454545454545
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -40.740740740740
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 16.129032258064
$ws.Range("I24").Value = 439
$ws.Range("J24").Value = 366
$ws.Range("K24").Value = 19.945355191256
$ws.Range("L24").Value = 50.342465753424
$ws.Range("M24").Value = 7.598039215686
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 30
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -2.325581395348
$ws.Range("I25").Value = 165
$ws.Range("J25").Value = 152
$ws.Range("K25").Value = 8.552631578947
$ws.Range("L25").Value = 25
$ws.Range("M25").Value = 35.245901639344
$ws.Range("D26").Value = 1
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = 150
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("E15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 166.666666666667
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = -11.764705882352
$ws.Range("L27").Value = -11.764705882352

$excel.CutCopyMode = $false
